$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '66.030.68'
$c.ClearFormats()
$ws.Range("E2").Value = '  -1.44%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.517.02'
$c.ClearFormats()
$ws.Range("E3").Value = '  -3.75%  '
$ws.Range("E4").Value = '  -0.03%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '578.21'
$c.ClearFormats()
$ws.Range("E5").Value = '  -2.68%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '166.74'
$c.ClearFormats()
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -1.74%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '2.519.24'
$c.ClearFormats()
$ws.Range("E9").Value = '  -3.72%  '
$ws.Range("E10").Value = '  -0.42%  '
$ws.Range("E11").Value = '  -0.15%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.349'
$c.ClearFormats()
$ws.Range("E12").Value = '  -3.73%  '
$ws.Range("E13").Value = '  -2.01%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '26.37'
$c.ClearFormats()
$ws.Range("E14").Value = '  -4.44%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.980.97'
$c.ClearFormats()
$ws.Range("E15").Value = '  -3.96%  '
$ws.Range("E16").Value = '  -3.39%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '65.870.28'
$c.ClearFormats()
$ws.Range("E17").Value = '  -2.12%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.512.33'
$c.ClearFormats()
$ws.Range("E18").Value = '  -3.87%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '11.29'
$c.ClearFormats()
$ws.Range("E19").Value = '  -5.74%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.65'
$c.ClearFormats()
$ws.Range("E20").Value = '  -4.14%  '
$ws.Range("E21").Value = '  -2.88%  '
$ws.Range("E22").Value = '  -2.69%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.56'
$c.ClearFormats()
$ws.Range("E23").Value = '  -2.01%  '
$ws.Range("E24").Value = '  +0.03%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.93'
$c.ClearFormats()
$ws.Range("E25").Value = '  +0.74%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '69.00'
$c.ClearFormats()
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("E27").Value = '  -3.03%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range("E28").Value = '  -0.29%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.644.70'
$c.ClearFormats()
$ws.Range("E29").Value = '  -4.05%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0₃0975'
$c.ClearFormats()
$ws.Range("E30").Value = '  -2.29%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '8.11'
$c.ClearFormats()
$ws.Range("E31").Value = '  +2.99%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '518.90'
$c.ClearFormats()
$ws.Range("E32").Value = '  -4.94%  '
$ws.Range("E33").Value = '  -2.61%  '
$ws.Range("E34").Value = '  -4.55%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.131'
$c.ClearFormats()
$ws.Range("E35").Value = '  -3.47%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.44'
$c.ClearFormats()
$ws.Range("E37").Value = '  -3.47%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '155.85'
$c.ClearFormats()
$ws.Range("E38").Value = '  -1.37%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '18.50'
$c.ClearFormats()
$ws.Range("E39").Value = '  -2.40%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '18.28'
$c.ClearFormats()
$ws.Range("E40").Value = '  +0.77%  '
$ws.Range("E41").Value = '  -3.21%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.75'
$c.ClearFormats()
$ws.Range("E42").Value = '  -2.26%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.04'
$c.ClearFormats()
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("E44").Value = '  +0.12%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.ClearFormats()
$ws.Range("E45").Value = '  +0.49%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '146.94'
$c.ClearFormats()
$ws.Range("E46").Value = '  -2.93%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0₆0280'
$c.ClearFormats()
$ws.Range("E47").Value = '  -5.65%  '
$ws.Range("E48").Value = '  -3.89%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '3.70'
$c.ClearFormats()
$ws.Range("E49").Value = '  -1.59%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.71'
$c.ClearFormats()
$ws.Range("E50").Value = '  +1.33%  '
$ws.Range("E51").Value = '  -2.19%  '
